$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (Marco): F9 "Incomplete" (red) -> "In Progress" (yellow) ---
# Capture the "In Progress" (yellow) formatting currently on F7 before we
# change F7's own formatting below.
$ws.Range("F7").Copy()
$ws.Range("F9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F9").Value = "In Progress"

# --- Row 7 (Ping): F7 "In Progress" (yellow) -> "Complete 0.3.1.4b" (green) ---
# Reuse the "Complete" (green) formatting already present on F3.
$ws.Range("F3").Copy()
$ws.Range("F7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F7").Value = "Complete 0.3.1.4b"

# New G7 line-complete number, matching the plain (unstyled) number cells
# used alongside other "Complete" rows (e.g. G3).
$ws.Range("G7").Value = 194

$excel.CutCopyMode = 0

# --- Selection update recorded in the sheet view ---
$ws.Range("G8").Select() | Out-Null
